$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-07-31T04:59:12.791Z", "arpitsin28@gmail.com", "8303884098", 582834),
    @("2025-07-31T05:00:19.028Z", "arpitsin28@gmail.com", "9473733115", 947559),
    @("2025-07-31T05:03:41.962Z", "arpitsin28@gmail.com", "9473733115", 140001),
    @("2025-07-31T05:06:01.226Z", "arpitsin28@gmail.com", "8303884098", 442704),
    @("2025-07-31T05:22:19.516Z", "arpitsin28@gmail.com", "9473733115", 269016),
    @("2025-07-31T05:36:50.240Z", "arpitsin28@gmail.com", "8303884098", 656820)
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    # Phone numbers are stored as text (like the existing rows), so force
    # the numeric-looking string to stay text with a leading quote prefix.
    $ws.Cells.Item($row, 3).Value = "'" + $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}
